$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 updates
$ws.Range("D15").Value = "image_20250807111728_ppp0.jpg"
# Leading apostrophe forces these numeric-looking strings to stay text,
# matching the source data (coords/confidence are stored as plain text).
$ws.Range("I15").Value = "'794,481,830,525"
$ws.Range("J15").Value = "'0.70"

# Row 16 updates
$ws.Range("D16").Value = "image_20250808221835_ppp0.jpg"
$ws.Range("I16").Value = "'1182,405,1231,455"
$ws.Range("J16").Value = "'0.76"
